$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.752.15'
$ws.Range("E2").Value = '  +0.55%  '
$ws.Range("D3").Value = '1.850.45'
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("D4").Value = '''1.032'
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").Value = '''322.10'
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("D6").Value = '''1.030'
$ws.Range("E6").Value = '  +0.17%  '
$ws.Range("D7").Value = '''0.4393'
$ws.Range("E7").Value = '  +0.40%  '
$ws.Range("D8").Value = '''0.3802'
$ws.Range("E8").Value = '  +1.49%  '
$ws.Range("D9").Value = '''0.07407'
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").Value = '''0.8846'
$ws.Range("E10").Value = '  +0.93%  '
$ws.Range("D11").Value = '''21.55'
$ws.Range("E11").Value = '  +0.16%  '
$ws.Range("D12").Value = '1.855.96'
$ws.Range("E12").Value = '  +0.22%  '
$ws.Range("D13").Value = '''5.509'
$ws.Range("E13").Value = '  +0.34%  '
$ws.Range("D14").Value = '''6.695'
$ws.Range("E14").Value = '  +0.09%  '
$ws.Range("D15").Value = '''0.07171'
$ws.Range("E15").Value = '  +0.26%  '
$ws.Range("D16").Value = '''85.11'
$ws.Range("E16").Value = '  +2.73%  '
$ws.Range("D17").Value = '''1.038'
$ws.Range("E17").Value = '  +0.40%  '
$ws.Range("D18").Value = '''0.000009080'
$ws.Range("E18").Value = '  +0.62%  '
$ws.Range("D19").Value = '''1.031'
$ws.Range("E19").Value = '  +0.42%  '
$ws.Range("D20").Value = '''15.49'
$ws.Range("E20").Value = '  +0.38%  '
$ws.Range("D21").Value = '27.783.37'
$ws.Range("E21").Value = '  +0.63%  '
$ws.Range("D22").Value = '''5.277'
$ws.Range("E22").Value = '  +0.40%  '
$ws.Range("E23").Value = '  +0.45%  '
$ws.Range("D24").Value = '2.086.57'
$ws.Range("E24").Value = '  +0.96%  '
$ws.Range("D25").Value = '''2.049'
$ws.Range("E25").Value = '  +6.23%  '
$ws.Range("D26").Value = '''158.47'
$ws.Range("E26").Value = '  +0.60%  '
$ws.Range("D27").Value = '''18.71'
$ws.Range("E27").Value = '  -0.16%  '
$ws.Range("D28").Value = '''1.995'
$ws.Range("E28").Value = '  +2.33%  '
$ws.Range("D29").Value = '''5.340'
$ws.Range("E29").Value = '  +1.40%  '
$ws.Range("D30").Value = '''117.81'
$ws.Range("E30").Value = '  +1.58%  '
$ws.Range("D31").Value = '''0.09070'
$ws.Range("E31").Value = '  -0.11%  '
$ws.Range("B32").Value = 'ARBITRUM'
$ws.Range("C32").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D32").Value = '''1.213'
$ws.Range("E32").Value = '  +0.52%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = '''0.7726'
$ws.Range("E33").Value = '  +0.68%  '
$ws.Range("D34").Value = '''3.010'
$ws.Range("E34").Value = '  +4.49%  '
$ws.Range("E35").Value = '  +1.31%  '
$ws.Range("D36").Value = '''1.032'
$ws.Range("E36").Value = '  +0.39%  '
$ws.Range("D37").Value = '''1.152'
$ws.Range("E37").Value = '  +0.50%  '
$ws.Range("D38").Value = '''0.01976'
$ws.Range("E38").Value = '  -0.32%  '
$ws.Range("D39").Value = '''0.05286'
$ws.Range("E39").Value = '  +0.19%  '
$ws.Range("D40").Value = '''2.857'
$ws.Range("E40").Value = '  +2.19%  '
$ws.Range("D41").Value = '''0.5183'
$ws.Range("D42").Value = '''0.1672'
$ws.Range("E42").Value = '  -0.09%  '
$ws.Range("D43").Value = '''6.869'
$ws.Range("E43").Value = '  +2.65%  '
$ws.Range("D44").Value = '''8.738'
$ws.Range("E44").Value = '  +2.04%  '
$ws.Range("D45").Value = '''110.41'
$ws.Range("E45").Value = '  +1.37%  '
$ws.Range("D46").Value = '''10.70'
$ws.Range("E46").Value = '  +1.31%  '
$ws.Range("D47").Value = '''1.033'
$ws.Range("E47").Value = '  +0.29%  '
$ws.Range("D48").Value = '''0.06573'
$ws.Range("E48").Value = '  +3.13%  '
$ws.Range("D49").Value = '''1.709'
$ws.Range("E49").Value = '  -0.40%  '
$ws.Range("D50").Value = '''0.4707'
$ws.Range("E50").Value = '  +1.18%  '
$ws.Range("D51").Value = '''1.895'
$ws.Range("E51").Value = '  +0.36%  '
